$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for P2, P4 and P5 (keep header, P3 and P6),
# working from bottom to top so row numbers don't shift underneath us.
$ws.Rows.Item(5).Delete()   # P5
$ws.Rows.Item(4).Delete()   # P4
$ws.Rows.Item(2).Delete()   # P2
